$wb = $excel.ActiveWorkbook

# Update the "展览" (Exhibition) sheet - F3:F6 "想去人数" (number of people interested) counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 289
$wsExhibit.Range("F4").Value = 2754
$wsExhibit.Range("F5").Value = 60
$wsExhibit.Range("F6").Value = 580

# Update the "全部类型" (All types) sheet - F5:F8 same counts for the same events
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 289
$wsAll.Range("F6").Value = 2754
$wsAll.Range("F7").Value = 60
$wsAll.Range("F8").Value = 580
